$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the event name in A4: "Shadow Assembly 6PM" -> "Shadow Assembly 7PM"
$ws.Range("A4").Value = "Shadow Assembly 7PM"

# Match the selection shown in the saved file (active cell A5)
$ws.Range("A5").Select()
